$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data per row: row number, new date text (dash-separated), D, E, F, G, H
$rows = @(
    @(3,  "28-07-2022", 1, 0, 0, 1, 1),
    @(4,  "01-08-2022", 1, 1, 0, 0, 0),
    @(5,  "04-08-2022", 1, 1, 0, 0, 0),
    @(6,  "08-08-2022", 0, 0, 0, 0, 1),
    @(7,  "11-08-2022", 0, 0, 0, 0, 1),
    @(8,  "15-08-2022", 0, 0, 0, 0, 1),
    @(9,  "18-08-2022", 0, 0, 0, 0, 1),
    @(10, "22-08-2022", 1, 1, 0, 0, 0),
    @(11, "25-08-2022", 0, 0, 0, 0, 1),
    @(12, "29-08-2022", 0, 0, 0, 0, 1),
    @(13, "01-09-2022", 0, 0, 0, 0, 1),
    @(14, "05-09-2022", 0, 0, 0, 0, 1),
    @(15, "08-09-2022", 0, 0, 0, 0, 1),
    @(16, "12-09-2022", 0, 0, 0, 0, 1),
    @(17, "15-09-2022", 0, 0, 0, 0, 1),
    @(18, "19-09-2022", 0, 0, 0, 0, 1),
    @(19, "22-09-2022", 0, 0, 0, 0, 1),
    @(20, "26-09-2022", 0, 0, 0, 0, 1),
    @(21, "29-09-2022", 0, 0, 0, 0, 1)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $dateText = $r[1]
    $dCell = $ws.Cells.Item($rowNum, 4)
    $eCell = $ws.Cells.Item($rowNum, 5)
    $fCell = $ws.Cells.Item($rowNum, 6)
    $gCell = $ws.Cells.Item($rowNum, 7)
    $hCell = $ws.Cells.Item($rowNum, 8)

    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $dateText

    $dCell.Value = $r[2]
    $eCell.Value = $r[3]
    $fCell.Value = $r[4]
    $gCell.Value = $r[5]
    $hCell.Value = $r[6]
}
